$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new log row (row 47) matching the existing data-row style (s="3":
# centered horizontal/vertical alignment, no border) by setting alignment
# explicitly on the whole row range before writing values.
$row = $ws.Range("A47:H47")
$row.HorizontalAlignment = -4108  # xlCenter
$row.VerticalAlignment = -4108    # xlCenter

$ws.Range("A47").Value = "2025-08-23 06:44:59 UTC"
$ws.Range("B47").Value = "2025-08-23 12:14:59 IST"
$ws.Range("C47").Value = "SKIPPED"
$ws.Range("D47").Value = "No change in PDF. Skipping download & Excel update."
$ws.Range("E47").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Range("G47").Value = 0
